# Update metricas_retencao_anual data to reflect refreshed BIBI metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 (cohort_year=2022, period_index=3): num_customers 47 -> 49,
# retention_rate recalculated as num_customers / cohort_size (2312).
$ws.Range("C31").Value = 49
$ws.Range("E31").Value = 0.02119377162629758

# Row 37 (cohort_year=2025, period_index=0): num_customers and
# cohort_size both updated from 851 -> 854 (retention_rate stays 1).
$ws.Range("C37").Value = 854
$ws.Range("D37").Value = 854
